# PSP_Sheet_정석준.xlsx - "Add files via upload" edit
# Fills in the two Time-Recording-Log rows (6 & 7) on the first sheet with
# real log entries (date / start / stop / interruption / delta / activity)
# and moves the sheet's active-cell selection to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 6 : 2019-09-06, 14:00 -> 15:40, 0 min interruption, 100 delta ---
$ws.Range("A6").Value = 43714
$ws.Range("B6").Value = 0.58333333333333337
$ws.Range("C6").Value = 0.65277777777777779
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = "아이디어 선정 및 프로젝트 계획서 작성"

# --- Row 7 : 2019-09-20, 13:00 -> 15:00, 0 min interruption, 120 delta ---
$ws.Range("A7").Value = 43728
$ws.Range("B7").Value = 0.54166666666666663
$ws.Range("C7").Value = 0.625
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 120
$ws.Range("F7").Value = "아이디어 정리 및 계획서 종합 수정"

# The activity note in F7 uses a mixed-font run: "아이디어" stays in the
# sheet's default font, the remainder (" 정리 및 계획서 종합 수정") is set
# in 돋움 10pt, matching the author's manual emphasis.
$note = $ws.Range("F7").Characters(5, 15)
$note.Font.Name = "돋움"
$note.Font.Size = 10

# Author left the cursor on B7 when saving.
$ws.Range("B7").Select()
